# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest scraped numbers (site regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet ---
$ws1.Range("F4").Value = 233
$ws1.Range("F5").Value = 1775
$ws1.Range("F8").Value = 485
$ws1.Range("F9").Value = 4507
$ws1.Range("F13").Value = 982
$ws1.Range("F14").Value = 1285
$ws1.Range("F17").Value = 2970
$ws1.Range("F18").Value = 1803
$ws1.Range("F21").Value = 168
$ws1.Range("F22").Value = 18
$ws1.Range("F24").Value = 926
$ws1.Range("F27").Value = 2321
$ws1.Range("F28").Value = 997
$ws1.Range("F29").Value = 2392
$ws1.Range("F31").Value = 1110
$ws1.Range("F32").Value = 568
$ws1.Range("F34").Value = 884
$ws1.Range("F35").Value = 415
$ws1.Range("F36").Value = 1102
$ws1.Range("F37").Value = 902
$ws1.Range("F38").Value = 1173
$ws1.Range("F40").Value = 841
$ws1.Range("F41").Value = 518
$ws1.Range("F42").Value = 358
$ws1.Range("F43").Value = 279
$ws1.Range("F44").Value = 3478

# --- 全部类型 sheet (same events, different row offsets) ---
$ws4.Range("F4").Value = 233
$ws4.Range("F6").Value = 1775
$ws4.Range("F9").Value = 485
$ws4.Range("F10").Value = 4507
$ws4.Range("F15").Value = 1285
$ws4.Range("F16").Value = 2970
$ws4.Range("F18").Value = 1803
$ws4.Range("F22").Value = 168
$ws4.Range("F25").Value = 18
$ws4.Range("F26").Value = 926
$ws4.Range("F28").Value = 2321
$ws4.Range("F31").Value = 997
$ws4.Range("F33").Value = 2392
$ws4.Range("F34").Value = 1110
$ws4.Range("F35").Value = 568
$ws4.Range("F36").Value = 884
$ws4.Range("F37").Value = 1102
$ws4.Range("F38").Value = 902
$ws4.Range("F40").Value = 1173
$ws4.Range("F41").Value = 841
$ws4.Range("F42").Value = 518
$ws4.Range("F44").Value = 358
$ws4.Range("F47").Value = 279
$ws4.Range("F48").Value = 3478

$wb.Save()
